$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1) Insert a new row at 25 (pushes CreateDate..LastUpdateEmpNo from 25-28 to 26-29)
$ws1.Rows.Item(25).Insert()

# 2) Fix the formula chain
$ws1.Range("A26").Formula = "=A25+1"
$ws1.Range("A25").Formula = "=A24+1"

# 3) Populate the new row 25 (GroupNo / 課組別)
$ws1.Range("B25").Value = "GroupNo"
$ws1.Range("C25").Value = "課組別"
$ws1.Range("D25").Value = "VARCHAR2"
$ws1.Range("E25").Value = 1
$ws1.Range("G25").Value = "CdBranchGroup：`nBranchNo單位別：0000`n放款管理課：1`n放款服務課：2`n放款推展課：3`n放款審查課：4`n投資資訊規劃課：5`n專案管理課：6`n軟體測試課：7"

Write-Host "values set"
